$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 16 entry for the Python tasks on lists/tuples/sets/dictionaries
# Copy cell A15 (date, with its date number-format/style) down to A16 first so
# the new date cell keeps the same style as the rows above it, then overwrite
# the value with the actual serial date for this entry.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 43325

$ws.Range("B16").Value = "python"
$ws.Range("C16").Value = "tasks on lists,tuples,sets,dictionaries,methods etc."

# Move the active selection to the newly added last cell, which also updates
# the sheet view (dimension grows automatically, scroll/topLeftCell resets).
$ws.Range("C16").Select()
